$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.833.37"
$ws.Range("E2").Value = "  -2.48%  "

# Row 3
$ws.Range("D3").Value = "3.363.75"
$ws.Range("E3").Value = "  +0.05%  "

# Row 4
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").Value = "'571.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.48%  "

# Row 6
$ws.Range("D6").Value = "'134.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.92%  "

# Row 7
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("D8").Value = "3.362.93"
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("D9").Value = "'0.477"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.47%  "

# Row 10
$ws.Range("D10").Value = "'7.61"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.42%  "

# Row 11
$ws.Range("E11").Value = "  +3.93%  "

# Row 12
$ws.Range("D12").Value = "'0.390"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.44%  "

# Row 13
$ws.Range("D13").Value = "3.934.54"
$ws.Range("E13").Value = "  -0.29%  "

# Row 15
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000172"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.93%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.362.07"
$ws.Range("E16").Value = "  -0.23%  "

# Row 17
$ws.Range("D17").Value = "'25.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.99%  "

# Row 18
$ws.Range("D18").Value = "60.930.54"
$ws.Range("E18").Value = "  -2.45%  "

# Row 19
$ws.Range("E19").Value = "  +7.05%  "

# Row 20
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'9.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.17%  "

# Row 21
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").Value = "'5.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.28%  "

# Row 22
$ws.Range("D22").Value = "'372.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.92%  "

# Row 23
$ws.Range("D23").Value = "'0.574"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.59%  "

# Row 24
$ws.Range("D24").Value = "3.497.18"
$ws.Range("E24").Value = "  -0.11%  "

# Row 25
$ws.Range("E25").Value = "  -0.03%  "

# Row 26
$ws.Range("D26").Value = "'70.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.76%  "

# Row 27
$ws.Range("D27").Value = "'0.0000117"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +13.51%  "

# Row 28
$ws.Range("E28").Value = "  +23.59%  "

# Row 29
$ws.Range("D29").Value = "'7.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +12.50%  "

# Row 30
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.06%  "

# Row 31
$ws.Range("D31").Value = "'8.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.67%  "

# Row 32
$ws.Range("E32").Value = "  +2.37%  "

# Row 33
$ws.Range("D33").Value = "'0.155"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.82%  "

# Row 34
$ws.Range("E34").Value = "  -0.04%  "

# Row 35
$ws.Range("D35").Value = "3.394.91"
$ws.Range("E35").Value = "  +0.01%  "

# Row 36
$ws.Range("D36").Value = "'23.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.72%  "

# Row 37
$ws.Range("D37").Value = "'5.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.69%  "

# Row 38
$ws.Range("E38").Value = "  +5.36%  "

# Row 39
$ws.Range("E39").Value = "  +6.20%  "

# Row 40
$ws.Range("D40").Value = "'162.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.26%  "

# Row 41
$ws.Range("D41").Value = "'0.0783"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.92%  "

# Row 42
$ws.Range("E42").Value = "  -0.11%  "

# Row 43
$ws.Range("E43").Value = "  +5.07%  "

# Row 44
$ws.Range("E44").Value = "  +12.26%  "

# Row 45
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'41.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.36%  "

# Row 46
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.757"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.52%  "

# Row 47
$ws.Range("D47").Value = "'1.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.98%  "

# Row 48
$ws.Range("D48").Value = "'23.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.58%  "

# Row 49
$ws.Range("D49").Value = "'6.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.46%  "

# Row 50
$ws.Range("D50").Value = "'23.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +17.64%  "

# Row 51
$ws.Range("E51").Value = "  +15.67%  "
